$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The title block (company name / "Trial Balance" / period) above the
# trial-balance table was left over from an older template and was no
# longer wanted on this page - clear it out.
$ws.Range("E2:E8").ClearContents()

# Column H held a redundant restatement of the account name (and a
# G56-F56 rounding-check formula in H56) that Excel was flagging with a
# warning on this client admin page. Select the whole column, as the
# user did in the workbook, and clear it.
$ws.Columns("H:H").Select()
$excel.Selection.ClearContents()
